$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the "Run the program." list paragraph entirely, merging
#    it away so the following paragraph ("On the computer screen,
#    press the "on" button. ") takes its place.
# ------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.ClearFormatting()
$found = $findRange.Find.Execute("Run the program.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find the 'Run the program.' paragraph."
}

$paraRunProgram = $findRange.Paragraphs(1)

# Delete the run's text together with its own trailing paragraph mark;
# this merges paragraph away and leaves the next paragraph's content
# sitting right where the bookmark needs to go.
$mergeRange = $d.Range($paraRunProgram.Range.Start, $paraRunProgram.Range.End)
$mergeRange.Delete()

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark here (right before the "On the
#    computer screen..." run). Re-adding a bookmark with a name that
#    already exists elsewhere relocates it, so the old one (around
#    "Turn" / " box on.") is automatically removed.
# ------------------------------------------------------------------
$bmRange = $d.Range($paraRunProgram.Range.Start, $paraRunProgram.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)
